# Update "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet, which mirrors the same events.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value = 7872
$wsExhibit.Range("F16").Value = 8492
$wsExhibit.Range("F19").Value = 1126
$wsExhibit.Range("F24").Value = 1034
$wsExhibit.Range("F25").Value = 9
$wsExhibit.Range("F26").Value = 555
$wsExhibit.Range("F29").Value = 571
$wsExhibit.Range("F33").Value = 103
$wsExhibit.Range("F38").Value = 3460
$wsExhibit.Range("F41").Value = 753
$wsExhibit.Range("F46").Value = 37

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F16").Value = 8492
$wsAll.Range("F18").Value = 1126
$wsAll.Range("F21").Value = 1034
$wsAll.Range("F22").Value = 555
$wsAll.Range("F24").Value = 571
$wsAll.Range("F27").Value = 103
$wsAll.Range("F32").Value = 3460
$wsAll.Range("F35").Value = 753
